$d = $word.ActiveDocument

# The template's second paragraph holds a complex field whose code is
# `{ m:'doc.html'.fromHTMLURI() }`, built from <w:fldChar>/<w:instrText>
# runs (with a _GoBack bookmark sitting between "doc.html" and
# "'.fromHTMLURI()"). The parser was switched to
# TokenIteratorFieldRewriterSplit, which expects the field markers to be
# rewritten as plain literal text ("{" ... "}") split across the same
# run boundaries, instead of an actual Word field.

$f = $d.Fields.Item(1)

# Deleting the Field removes every run that belongs to it (fldChar
# begin/end, all instrText runs and the bookmark) but keeps the now
# empty host paragraph (and its trailing paragraph mark) intact. The
# field's opening <w:fldChar begin> sits one position before its code.
$start = $f.Code.Start - 1
$f.Delete()

# Re-insert the field's former content as plain text, using the same
# literal punctuation the M2Doc template syntax expects: a leading "{",
# the code itself and a trailing "}".
$full = "{m:'doc.html'.fromHTMLURI()}"
$ins = $d.Range($start, $start)
$ins.InsertAfter($full)

# Word merges adjacent same-formatted text into a single run, but the
# target markup keeps each former field-code token in its own <w:r>.
# Dropping (and immediately removing) a temporary bookmark at each
# boundary forces the run to split there without altering the visible
# text or formatting. The boundary between "doc.html" and
# "'.fromHTMLURI()" keeps its bookmark permanently: that's where the
# original _GoBack bookmark used to live.
$boundaries = @(1, 2, 3, 4)
$i = 0
foreach ($b in $boundaries) {
    $i = $i + 1
    $pos = $start + $b
    $name = "M2DocSplit" + $i
    $d.Bookmarks.Add($name, $d.Range($pos, $pos))
}

$d.Bookmarks.Add("_GoBack", $d.Range($start + 12, $start + 12))

$lastName = "M2DocSplitLast"
$d.Bookmarks.Add($lastName, $d.Range($start + 27, $start + 27))

$d.Bookmarks.ShowHidden = $true
for ($j = 1; $j -le $i; $j++) {
    $d.Bookmarks.Item("M2DocSplit" + $j).Delete()
}
$d.Bookmarks.Item($lastName).Delete()
